$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: for each destination row, the source row whose original
# (D, J, K, L, M, O, P) values should be written into it.
$rowMap = @{
    2 = 10
    3 = 17
    4 = 4
    5 = 13
    6 = 29
    7 = 15
    8 = 32
    9 = 37
    10 = 31
    11 = 35
    12 = 46
    13 = 22
    14 = 50
    15 = 19
    16 = 28
    17 = 42
    18 = 7
    19 = 5
    20 = 45
    21 = 21
    22 = 36
    23 = 9
    24 = 20
    25 = 25
    26 = 18
    27 = 26
    28 = 49
    29 = 11
    30 = 40
    31 = 43
    32 = 8
    33 = 41
    34 = 30
    35 = 47
    36 = 34
    37 = 23
    38 = 12
    39 = 38
    40 = 3
    41 = 27
    42 = 16
    43 = 48
    44 = 14
    45 = 33
    46 = 39
    47 = 44
    48 = 6
    49 = 24
    50 = 2
}

# Columns whose values get shuffled between rows.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot original values before any writes, so later writes don't
# clobber values still needed as a source for other rows.
$snapshot = @{}
for ($r = 2; $r -le 50; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
